# Update the LR-pair sheet with newly computed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Il1a-Il1r1, Target cluster ECs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02286966666666667
$ws.Range("H2").Value = 0.068609
$ws.Range("M2").Value = 6.382924
$ws.Range("N2").Value = 19.148772
$ws.Range("O2").Value = 0.1363153751023214
$ws.Range("P2").Value = 0.1363153751023214
$ws.Range("Q2").Value = 0.1459753442386667
$ws.Range("R2").Value = 1.313778098148
$ws.Range("S2").Value = 0.1363153751023214
$ws.Range("T2").Value = 0.1363153751023214

# Row 3 (Il1a-Il1r1, Target cluster FAPs)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02286966666666667
$ws.Range("H3").Value = 0.068609
$ws.Range("O3").Value = 0.6265841681043937
$ws.Range("P3").Value = 0.6265841681043938
$ws.Range("Q3").Value = 0.6709869636119999
$ws.Range("R3").Value = 6.038882672508
$ws.Range("S3").Value = 0.6265841681043937
$ws.Range("T3").Value = 0.6265841681043938

# Row 4 (Il1a-Il1r1, Target cluster MuSCs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02286966666666667
$ws.Range("H4").Value = 0.068609
$ws.Range("O4").Value = 0.2371004567932849
$ws.Range("P4").Value = 0.2371004567932849
$ws.Range("Q4").Value = 0.2539025460155556
$ws.Range("R4").Value = 2.28512291414
$ws.Range("S4").Value = 0.2371004567932849
$ws.Range("T4").Value = 0.2371004567932849
